$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.044.45'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '3.114.26'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.56%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '3.109.40'
$ws.Range("E8").Value = '  +0.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.521'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.45'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.37%  '
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.479'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000246'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.16'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.124'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.92%  '
$ws.Range("E16").Value = '  +1.17%  '
$ws.Range("D17").Value = '66.972.96'
$ws.Range("E17").Value = '  +0.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.11'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.88%  '
$ws.Range("D19").Value = '3.110.19'
$ws.Range("E19").Value = '  +1.13%  '
$ws.Range("E20").Value = '  +1.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '477.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.713'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.94'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.15%  '
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("D34").Value = '0.0₃0945'
$ws.Range("E34").Value = '  -7.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.981'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '47.72'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.09'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.97'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.61'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.75%  '
$ws.Range("D44").Value = '2.796.70'
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0355'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '379.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -11.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '136.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.21'
$ws.Range("D51").Style = "Normal"
